$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the time-range labels in column C for rows 2, 3, 6, 7
$ws.Range("C2").Value = "9:30-9:35"
$ws.Range("C3").Value = "9:35-9:40"
$ws.Range("C6").Value = "22:30-22:35"
$ws.Range("C7").Value = "22:35-22:40"

# Update the selection to C12 on the active sheet
$ws.Range("C12").Select()
